$wb = $excel.ActiveWorkbook

# --- "HRIS Seq Diagram" sheet (4th sheet) -------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws5 = $wb.Worksheets.Item(5)

# New "docker-compose 2" block of rows (order matters: it controls the
# order new shared strings are appended in).
$ws4.Range("A12").Value = "5.2.0"
$ws4.Range("B12").Value = "Declare all necesarry branch queues"
$ws4.Range("B13").Value = "CORP:agm_branch -> CORP:agm_branch"
$ws4.Range("C13").Value = "agm_branch_boot"
$ws4.Range("B21").Value = "Done 5.2.1"

# New "Status" column (D) on the existing rows.
$ws4.Range("D3").Value = "DONE"
$ws4.Range("D4").Value = "DONE"
$ws4.Range("D7").Value = "DONE"
$ws4.Range("D8").Value = "DONE"
$ws4.Range("D1").Value = "Status"

# Widen column C and give column E a width (closest the engine's pixel
# grid allows to the authored 42.1640625 / 23.6640625 character widths).
$ws4.Columns.Item(3).ColumnWidth = 41.25
$ws4.Columns.Item(5).ColumnWidth = 22.75

# Update the sheet's selection / scroll anchor: drop the old
# topLeftCell="B1" scroll (by re-selecting from the top) and move the
# active cell to D9.
$ws4.Range("D9").Select() | Out-Null

# --- "HRIS DB" sheet (5th / last sheet) ---------------------------------
# Clear its scroll anchor (topLeftCell="A2") without disturbing which
# sheet/cell is actually active in the workbook.
$ws5.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# Restore "HRIS Seq Diagram" as the active tab.
$ws4.Activate() | Out-Null
